$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "브리츠 BZ-505 Roll"
$ws.Range("C3").Value = "https://search.shopping.naver.com/gate.nhn?id=20665125832"
$ws.Range("D3").Value = "https://shopping-phinf.pstatic.net/main_2066512/20665125832.20190819141404.jpg"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "19900"
$ws.Range("I3").Value = "브리츠"
$ws.Range("J3").Value = "브리츠"

$ws.Range("B4").Value = "에이튜브 A450V 진공관 사운드바 스피커"
$ws.Range("C4").Value = "https://search.shopping.naver.com/gate.nhn?id=25429469494"
$ws.Range("D4").Value = "https://shopping-phinf.pstatic.net/main_2542946/25429469494.20210218191500.jpg"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "39000"
$ws.Range("I4").Value = "에이튜브"
$ws.Range("J4").Value = ""

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "19800"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "14550"

$ws.Range("B28").Value = "에프엔디 에프앤디 R27BT"
$ws.Range("C28").Value = "https://search.shopping.naver.com/gate.nhn?id=23013561490"
$ws.Range("D28").Value = "https://shopping-phinf.pstatic.net/main_2301356/23013561490.20210528134443.jpg"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "79800"
$ws.Range("I28").Value = "에프앤디"
$ws.Range("J28").Value = "에프엔디"

$ws.Range("B29").Value = "브리츠 BR-1000A2"
$ws.Range("C29").Value = "https://search.shopping.naver.com/gate.nhn?id=7872712963"
$ws.Range("D29").Value = "https://shopping-phinf.pstatic.net/main_7872712/7872712963.20140925112335.jpg"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "75000"
$ws.Range("I29").Value = "브리츠"
$ws.Range("J29").Value = "브리츠"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "18580"

$ws.Range("B36").Value = "브리츠 Z2100BT Edition"
$ws.Range("C36").Value = "https://search.shopping.naver.com/gate.nhn?id=24727397522"
$ws.Range("D36").Value = "https://shopping-phinf.pstatic.net/main_2472739/24727397522.20201105145232.jpg"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "37900"
$ws.Range("I36").Value = "브리츠"
$ws.Range("J36").Value = "브리츠"

$ws.Range("B37").Value = "앱코 SP200 Hybrid"
$ws.Range("C37").Value = "https://search.shopping.naver.com/gate.nhn?id=21894657996"
$ws.Range("D37").Value = "https://shopping-phinf.pstatic.net/main_2189465/21894657996.20200303105232.jpg"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "11000"
$ws.Range("I37").Value = "앱코"
$ws.Range("J37").Value = "앱코"

$ws.Range("B38").Value = "맥스틸 SB-100"
$ws.Range("C38").Value = "https://search.shopping.naver.com/gate.nhn?id=9856927406"
$ws.Range("D38").Value = "https://shopping-phinf.pstatic.net/main_9856927/9856927406.20210506155134.jpg"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "14300"
$ws.Range("I38").Value = "맥스틸"
$ws.Range("J38").Value = "맥스틸"

$ws.Range("B39").Value = "오디오엔진 A2+ Wireless"
$ws.Range("C39").Value = "https://search.shopping.naver.com/gate.nhn?id=18591127164"
$ws.Range("D39").Value = "https://shopping-phinf.pstatic.net/main_1859112/18591127164.20190416151115.jpg"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "338000"
$ws.Range("I39").Value = "오디오엔진"
$ws.Range("J39").Value = "오디오엔진"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "541090"

$ws.Range("B50").Value = "캔스톤 NX201 BOSS"
$ws.Range("C50").Value = "https://search.shopping.naver.com/gate.nhn?id=17169784911"
$ws.Range("D50").Value = "https://shopping-phinf.pstatic.net/main_1716978/17169784911.20201231114206.jpg"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "128000"
$ws.Range("I50").Value = "캔스톤"
$ws.Range("J50").Value = "캔스톤"

$ws.Range("B51").Value = "크리에이티브 GIGAWORKS T20 II"
$ws.Range("C51").Value = "https://search.shopping.naver.com/gate.nhn?id=5509585915"
$ws.Range("D51").Value = "https://shopping-phinf.pstatic.net/main_5509585/5509585915.20190924144502.jpg"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "109000"

$ws.Range("B52").Value = "크리에이티브 GIGAWORKS T40 II"
$ws.Range("C52").Value = "https://search.shopping.naver.com/gate.nhn?id=5509585925"
$ws.Range("D52").Value = "https://shopping-phinf.pstatic.net/main_5509585/5509585925.20190924144428.jpg"
$ws.Range("E52").NumberFormat = "@"
$ws.Range("E52").Value = "159000"
$ws.Range("I52").Value = "크리에이티브"
$ws.Range("J52").Value = "크리에이티브"

$ws.Range("B67").Value = "캔스톤 ATP-1"
$ws.Range("C67").Value = "https://search.shopping.naver.com/gate.nhn?id=20665881242"
$ws.Range("D67").Value = "https://shopping-phinf.pstatic.net/main_2066588/20665881242.20190819153106.jpg"
$ws.Range("E67").NumberFormat = "@"
$ws.Range("E67").Value = "36900"
$ws.Range("I67").Value = "캔스톤"
$ws.Range("J67").Value = "캔스톤"

$ws.Range("B68").Value = "몬스타기어 LIVE100 사운드바 usb 스피커"
$ws.Range("C68").Value = "https://search.shopping.naver.com/gate.nhn?id=19030267059"
$ws.Range("D68").Value = "https://shopping-phinf.pstatic.net/main_1903026/19030267059.20210309133924.jpg"
$ws.Range("E68").NumberFormat = "@"
$ws.Range("E68").Value = "17000"
$ws.Range("I68").Value = "몬스타기어"
$ws.Range("J68").Value = "몬스타"

$ws.Range("B72").Value = "인켈 IK-KS500"
$ws.Range("C72").Value = "https://search.shopping.naver.com/gate.nhn?id=21901171642"
$ws.Range("D72").Value = "https://shopping-phinf.pstatic.net/main_2190117/21901171642.20200313160748.jpg"
$ws.Range("E72").NumberFormat = "@"
$ws.Range("E72").Value = "9900"
$ws.Range("I72").Value = "인켈"
$ws.Range("J72").Value = "인켈"

$ws.Range("B73").Value = "브리츠 BE-100 Soundbar Plus"
$ws.Range("C73").Value = "https://search.shopping.naver.com/gate.nhn?id=7131613034"
$ws.Range("D73").Value = "https://shopping-phinf.pstatic.net/main_7131613/7131613034.20210928113139.jpg"
$ws.Range("E73").NumberFormat = "@"
$ws.Range("E73").Value = "26900"
$ws.Range("I73").Value = "브리츠"
$ws.Range("J73").Value = "브리츠"

$ws.Range("B74").Value = "몬스타기어 가츠 라이브200 블루투스 사운드바 스피커"
$ws.Range("C74").Value = "https://search.shopping.naver.com/gate.nhn?id=21952153318"
$ws.Range("D74").Value = "https://shopping-phinf.pstatic.net/main_2195215/21952153318.20201210183159.jpg"
$ws.Range("E74").NumberFormat = "@"
$ws.Range("E74").Value = "23500"
$ws.Range("I74").Value = "몬스타기어"
$ws.Range("J74").Value = "몬스타"

$ws.Range("B83").Value = "앱코 S1300"
$ws.Range("C83").Value = "https://search.shopping.naver.com/gate.nhn?id=24513037522"
$ws.Range("D83").Value = "https://shopping-phinf.pstatic.net/main_2451303/24513037522.20201019133951.jpg"
$ws.Range("E83").NumberFormat = "@"
$ws.Range("E83").Value = "21900"
$ws.Range("I83").Value = "앱코"
$ws.Range("J83").Value = "앱코"

$ws.Range("B84").Value = "컴소닉 PILLAR CS-50U"
$ws.Range("C84").Value = "https://search.shopping.naver.com/gate.nhn?id=18340521348"
$ws.Range("D84").Value = "https://shopping-phinf.pstatic.net/main_1834052/18340521348.20190327095751.jpg"
$ws.Range("E84").NumberFormat = "@"
$ws.Range("E84").Value = "15900"
$ws.Range("I84").Value = "PILLAR"
$ws.Range("J84").Value = "컴소닉"

$ws.Range("B99").Value = "아이리버 IR-S50 WOODEN BLOCK"
$ws.Range("C99").Value = "https://search.shopping.naver.com/gate.nhn?id=12442687510"
$ws.Range("D99").Value = "https://shopping-phinf.pstatic.net/main_1244268/12442687510.20181002114729.jpg"
$ws.Range("E99").NumberFormat = "@"
$ws.Range("E99").Value = "24900"
$ws.Range("I99").Value = "아이리버"
$ws.Range("J99").Value = "아이리버"

$ws.Range("B100").Value = "다름인터내셔널 디알고 BT-RGB2CH"
$ws.Range("C100").Value = "https://search.shopping.naver.com/gate.nhn?id=29893098619"
$ws.Range("D100").Value = "https://shopping-phinf.pstatic.net/main_2989309/29893098619.20211130101424.jpg"
$ws.Range("E100").NumberFormat = "@"
$ws.Range("E100").Value = "29360"
$ws.Range("I100").Value = "디알고"
$ws.Range("J100").Value = "다름인터내셔널"
